$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all changed Price (D) and Volume (E) cells as text to preserve exact formatting
$cells = @{
    'D2' = '34.196.95'
    'E2' = '  +0.32%  '
    'D3' = '1.788.41'
    'E3' = '  -0.10%  '
    'E5' = '  -0.42%  '
    'E6' = '  +0.52%  '
    'E7' = '  +0.18%  '
    'D8' = '32.29'
    'E8' = '  +0.20%  '
    'D9' = '0.295'
    'E9' = '  +0.00%  '
    'E10' = '  +0.23%  '
    'E11' = '  +0.71%  '
    'D12' = '2.046.93'
    'E12' = '  -0.02%  '
    'D13' = '11.15'
    'E13' = '  -1.83%  '
    'D14' = '1.786.40'
    'E14' = '  -0.17%  '
    'E15' = '  +0.43%  '
    'D16' = '34.177.76'
    'E16' = '  +0.30%  '
    'E17' = '  +0.34%  '
    'D18' = '67.97'
    'E18' = '  -0.13%  '
    'D19' = '0.0₃0804'
    'E19' = '  +2.86%  '
    'D20' = '245.95'
    'E20' = '  +0.86%  '
    'D21' = '11.02'
    'E21' = '  +0.85%  '
    'E22' = '  +0.17%  '
    'E23' = '  +1.81%  '
    'E24' = '  +0.54%  '
    'D25' = '161.94'
    'E25' = '  -0.01%  '
    'E26' = '  -0.44%  '
    'D27' = '16.32'
    'E27' = '  +0.24%  '
    'E28' = '  +0.80%  '
    'E29' = '  +0.30%  '
    'E31' = '  -0.23%  '
    'E32' = '  +2.93%  '
    'E33' = '  +3.95%  '
    'E34' = '  -1.58%  '
    'D35' = '1.442.19'
    'E35' = '  +2.06%  '
    'D36' = '2.57'
    'E36' = '  +8.41%  '
    'D37' = '0.667'
    'E37' = '  +2.95%  '
    'E38' = '  +1.04%  '
    'D39' = '0.0190'
    'E39' = '  -0.15%  '
    'D40' = '82.16'
    'E40' = '  +1.80%  '
    'E41' = '  +1.53%  '
    'E42' = '  +0.29%  '
    'E43' = '  +1.12%  '
    'D44' = '13.80'
    'E44' = '  +3.32%  '
    'D45' = '0.0520'
    'E45' = '  +2.47%  '
    'D46' = '6.11'
    'E46' = '  +1.02%  '
    'E47' = '  +0.71%  '
    'D48' = '1.946.76'
    'E48' = '  -0.05%  '
    'D49' = '105.32'
    'E49' = '  -1.63%  '
    'E50' = '  +0.21%  '
    'E51' = '  -6.98%  '
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
    $rng.Style = "Normal"
}
